$d = $word.ActiveDocument

# The last paragraph in the document body is an (almost) empty heading
# paragraph (numbered list style "a5") that only contains the _GoBack
# bookmark. Grab it and put the cursor at its very start so new text is
# typed before the bookmark.
$count = $d.Paragraphs.Count
$headingPara = $d.Paragraphs.Item($count)
$headingRange = $headingPara.Range
$headingRange.Collapse(1)  # wdCollapseStart
$headingRange.InsertBefore("Система")

# Insert a brand-new paragraph right after the heading paragraph for the
# description text. InsertParagraphAfter() adds a paragraph break after the
# heading paragraph's range and returns a Range positioned in the new
# paragraph.
$headingParaRange = $headingPara.Range
$headingParaRange.Collapse(0)  # wdCollapseEnd
$headingParaRange.InsertParagraphAfter()

# The newly created paragraph is now the second-to-last paragraph (the
# bookmark paragraph remains last). Reset its formatting: no list style,
# first-line indent 0.
$newCount = $d.Paragraphs.Count
$descPara = $d.Paragraphs.Item($newCount - 1)
$descPara.Style = $d.Styles.Item("Normal")
$descPara.Range.ParagraphFormat.FirstLineIndent = 0

$descRange = $descPara.Range
$descRange.Collapse(0)  # wdCollapseEnd, right before the paragraph mark
$descRange.InsertBefore("Распределенная система, ориентированная на управление задачами и проверку их решений в сфере олимпиадного программирования.")
